$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextCell "D2" "29.639.83"
Set-TextCell "E2" "  +5.66%  "
Set-TextCell "D3" "1.920.69"
Set-TextCell "E3" "  +3.96%  "
Set-TextCell "D4" "1.001"
Set-TextCell "E4" "  -0.44%  "
Set-TextCell "D5" "335.48"
Set-TextCell "E5" "  +1.30%  "
Set-TextCell "D6" "1.001"
Set-TextCell "E6" "  -0.62%  "
Set-TextCell "D7" "0.4675"
Set-TextCell "E7" "  +3.19%  "
Set-TextCell "E8" "  +6.38%  "
Set-TextCell "E9" "  +1.32%  "
Set-TextCell "D10" "0.08044"
Set-TextCell "E10" "  +4.13%  "
Set-TextCell "E11" "  +4.16%  "
Set-TextCell "D12" "22.47"
Set-TextCell "E12" "  +6.02%  "
Set-TextCell "D13" "1.915.37"
Set-TextCell "E13" "  +3.72%  "
Set-TextCell "D14" "6.007"
Set-TextCell "E14" "  +4.38%  "
Set-TextCell "D15" "7.190"
Set-TextCell "E15" "  +3.16%  "
Set-TextCell "D16" "89.92"
Set-TextCell "E16" "  +3.52%  "
Set-TextCell "D17" "1.001"
Set-TextCell "E17" "  -0.31%  "
Set-TextCell "D18" "0.00001037"
Set-TextCell "E18" "  +2.14%  "
Set-TextCell "D19" "0.06590"
Set-TextCell "E19" "  +0.77%  "
Set-TextCell "E20" "  +5.83%  "
Set-TextCell "D21" "1.002"
Set-TextCell "E21" "  -1.34%  "
Set-TextCell "D22" "29.611.49"
Set-TextCell "E22" "  +5.66%  "
Set-TextCell "D23" "5.576"
Set-TextCell "E23" "  +5.70%  "
Set-TextCell "D24" "11.66"
Set-TextCell "E24" "  +10.47%  "
Set-TextCell "E25" "  -2.49%  "
Set-TextCell "D26" "2.172.66"
Set-TextCell "E26" "  +5.13%  "
Set-TextCell "D27" "156.28"
Set-TextCell "E27" "  +0.17%  "
Set-TextCell "E28" "  +4.48%  "
Set-TextCell "D29" "2.147"
Set-TextCell "E29" "  +6.05%  "
Set-TextCell "D30" "5.755"
Set-TextCell "E30" "  +10.16%  "
Set-TextCell "D31" "117.49"
Set-TextCell "E31" "  +1.17%  "
Set-TextCell "D32" "1.076"
Set-TextCell "E32" "  +15.81%  "
Set-TextCell "D33" "0.09493"
Set-TextCell "E33" "  +2.72%  "
Set-TextCell "D34" "1.434"
Set-TextCell "E34" "  +5.20%  "
Set-TextCell "D35" "5.426"
Set-TextCell "E35" "  +5.51%  "
Set-TextCell "D36" "3.529"
Set-TextCell "E36" "  -2.41%  "
Set-TextCell "D37" "0.06146"
Set-TextCell "E37" "  +2.72%  "
Set-TextCell "D38" "0.02276"
Set-TextCell "E38" "  +4.41%  "
Set-TextCell "D39" "8.442"
Set-TextCell "E39" "  +4.16%  "
Set-TextCell "D40" "1.183"
Set-TextCell "E40" "  +2.64%  "
Set-TextCell "D41" "0.5901"
Set-TextCell "E41" "  +4.72%  "
Set-TextCell "D42" "0.1849"
Set-TextCell "E42" "  +3.85%  "
Set-TextCell "D43" "10.24"
Set-TextCell "E43" "  +3.62%  "
Set-TextCell "B44" "WEMIXTOKEN"
Set-TextCell "C44" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell "D44" "1.260"
Set-TextCell "E44" "  +1.45%  "
Set-TextCell "B45" "RenderToken"
Set-TextCell "C45" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell "D45" "2.353"
Set-TextCell "E45" "  +3.64%  "
Set-TextCell "D46" "0.07509"
Set-TextCell "E46" "  +4.82%  "
Set-TextCell "D47" "0.5585"
Set-TextCell "E47" "  +4.70%  "
Set-TextCell "D48" "12.21"
Set-TextCell "E48" "  +4.57%  "
Set-TextCell "E49" "  +3.88%  "
Set-TextCell "D50" "113.21"
Set-TextCell "E50" "  +3.47%  "
Set-TextCell "D51" "0.2997"
Set-TextCell "E51" "  +13.90%  "
